$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-sorts several data rows:
#   - Row 34 and Row 35 swap places entirely.
#   - Rows 37, 38, 39 rotate: new row37 = old row39, new row38 = old row37,
#     new row39 = old row38.
#
# We stage each source row in a scratch row far below the used range and
# copy it back in, which preserves each cell's original type (numbers stay
# numbers, date-look-alike text like "2025-06-26" stays text, etc.) instead
# of re-evaluating string contents.
#
# Range.Copy only overwrites cells that actually exist in the source range;
# cells that are blank in the source are left untouched in the destination.
# Because of that, after the moves below we explicitly blank out the one or
# two destination cells that must end up empty (they held a value in the
# row that used to occupy that slot, but the incoming row never had a value
# there).

$scratch = $ws.Range("A1000:AY1000")

# --- Swap rows 34 and 35 -------------------------------------------------
$ws.Range("A34:AY34").Copy($scratch)
$ws.Range("A35:AY35").Copy($ws.Range("A34:AY34"))
$scratch.Copy($ws.Range("A35:AY35"))
$scratch.Clear()

# --- Rotate rows 37 -> 38 -> 39 -> 37 ------------------------------------
# (new37 = old39, new38 = old37, new39 = old38)
$ws.Range("A39:AY39").Copy($scratch)
$ws.Range("A38:AY38").Copy($ws.Range("A39:AY39"))
$ws.Range("A37:AY37").Copy($ws.Range("A38:AY38"))
$scratch.Copy($ws.Range("A37:AY37"))
$scratch.Clear()

# New row 38 came from old row37, which had nothing in J (Enhet); old row38
# did, so clear the stale leftover.
$ws.Range("J38").Clear()
# New row 39 came from old row38, which had nothing in AC (Publik kommentar);
# old row39 did, so clear the stale leftover.
$ws.Range("AC39").Clear()
